$wb = $excel.ActiveWorkbook

# GameConfig is the second sheet in the workbook
$ws = $wb.Worksheets.Item("GameConfig")

# Fill in a new config row (row 7) mirroring the pattern of row 6
$ws.Range("A7").Value = "d13_talk_to_zora"
$ws.Range("B7").Value = "bool"
$ws.Range("D7").Value = $false

# Update the active selection on the GameConfig sheet
$ws.Activate()
$ws.Range("B9").Select()
